# Apply updated vm_pu results for the 380 kV case (commit: "case with 380 kV done")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout used by the data block: B,C,D,E,F,I,J,K,L,M,N (columns 2-6,9-14)
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

# New values for rows 2..25, one sub-array per row, matching $cols order
$data = @(
    @(1.02,1.021825799315812,1.02677758826004,1.025446449336963,1.020263360083504,1.028636915750121,1.027015082493927,1.029599243765538,1.028271995973042,1.023104162991646,1.028473561711347),
    @(1.02,1.023060463188128,1.027705593567972,1.026628730790911,1.022141080287777,1.028897005054352,1.027886049338967,1.030334508988781,1.029260564439648,1.02478514686409,1.029345765429206),
    @(1.02,1.023858039731887,1.028304458515929,1.027392797558433,1.023354415912725,1.0290628047428,1.028447841553504,1.030808041505883,1.029898725652598,1.025870786334476,1.029908355453126),
    @(1.02,1.024193026489813,1.02855583703917,1.027713788236555,1.023864114188115,1.029131910553565,1.028683595978233,1.031006582556831,1.030166651159147,1.026326705254306,1.03014444467618),
    @(1.02,1.024249253942322,1.028598022141548,1.027767671024671,1.023949672392328,1.029143478781323,1.028723155478028,1.031039887350632,1.030211616175844,1.026403228003625,1.030184060355002),
    @(1.02,1.023862517070864,1.028307818954823,1.027397087522814,1.023361228034744,1.029063730480139,1.02845099337359,1.030810696508415,1.029902307086387,1.025876880232675,1.029911511749158),
    @(1.02,1.022243338790567,1.027091548352016,1.025846205712676,1.020898299041431,1.028725331490365,1.027309800815086,1.029848194089711,1.02860640194622,1.0236726948276,1.028768698566325),
    @(1.02,1.019379708402245,1.024935820383436,1.023105906916156,1.01654489301558,1.02810986875443,1.025285071067572,1.028134918464202,1.026311128683648,1.019772269943535,1.02674109347033),
    @(1.02,1.017463299696199,1.023490085329444,1.021273774141786,1.013632788682517,1.027686607995972,1.023925747036516,1.026980983986623,1.024772826293016,1.017160253517119,1.025379839043232),
    @(1.02,1.016631667316579,1.022861989800422,1.02047913222983,1.012369300417548,1.027500242958489,1.023334839616875,1.026478491381429,1.024104744283474,1.016026275922723,1.024788092467278),
    @(1.02,1.016322483263182,1.022628370626005,1.020183763885027,1.011899589041938,1.027430552888163,1.023114998577414,1.026291414375148,1.023856285593376,1.015604607766131,1.024567939228324),
    @(1.02,1.016388817005393,1.022678497114309,1.02024713066587,1.01200036188231,1.02744552274765,1.023162171215063,1.026331562492417,1.023909594627311,1.015695077987818,1.02461517885653),
    @(1.02,1.016606115788449,1.022842685244841,1.020454721160823,1.012330482090261,1.027494491866494,1.023316674688986,1.026463036307154,1.024084212825376,1.015991430134346,1.024769901743108),
    @(1.02,1.016739963618289,1.022943804918448,1.020582597465959,1.012533827152888,1.027524601588488,1.023411822637842,1.026543984727426,1.024191760528374,1.016173961307901,1.024865184812964),
    @(1.02,1.017518453643074,1.023531725772547,1.021326483732828,1.013716587633471,1.027698911191676,1.023964914540589,1.027014272814801,1.024817122427375,1.017235448369582,1.025419062169654),
    @(1.02,1.018006289708512,1.023899952847401,1.021792747866732,1.01445781414853,1.027807422572484,1.024311232565399,1.027308511615725,1.025208860116064,1.017900489511488,1.025765872005782),
    @(1.02,1.0182906615504,1.024114532767528,1.02206458519537,1.014889916205605,1.02787041735638,1.024513011324873,1.027479863184576,1.025437162698913,1.018288112385699,1.025967937314239),
    @(1.02,1.018387595554256,1.024187665040885,1.022157253373022,1.015037211128079,1.027891846437184,1.024581775030929,1.027538243449922,1.02551497566112,1.018420233958563,1.026036798672646),
    @(1.02,1.017953967629386,1.023860466354108,1.021742735277013,1.014378312792265,1.027795811170348,1.024274098985875,1.027276970835412,1.025166850200882,1.017829166409295,1.025728685692316),
    @(1.02,1.016542134457368,1.0227943446998,1.020393596578721,1.012233280964485,1.027480084552064,1.023271187012147,1.026424332413117,1.024032800529045,1.015904174581177,1.02472434946855),
    @(1.02,1.015652843021633,1.022122198529749,1.019544162670643,1.010882317701868,1.02727887896836,1.022638579281378,1.025885761665303,1.023318021372023,1.014691196649576,1.024090843362242),
    @(1.02,1.016124428460452,1.022478690918667,1.019994576958136,1.011598711891996,1.02738579788973,1.022974131220721,1.026171504761201,1.023697107488535,1.01533447599965,1.024426871823827),
    @(1.02,1.017977610264047,1.023878309220668,1.021765334204512,1.01441423678178,1.027801058784053,1.024290878739784,1.027291223599491,1.025185833264249,1.017861395126363,1.025745489275401),
    @(1.02,1.02012129322304,1.025494627471248,1.023815248927243,1.017672025732556,1.028271258337103,1.025810170831217,1.028579899754605,1.026905925433863,1.020782637572594,1.027266938935868)
)

$rowStart = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $rowStart + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}

Write-Host "Updated $($data.Count) rows x $($cols.Count) cols"
